# Fruta / hortaliza, semanal
# A new weekly price record for Puerro (Vega Central Mapocho de Santiago)
# is inserted as row 46; every existing record from row 46 down to row 66
# shifts down one row (to 47..67), and the sheet's used range grows to
# A1:R67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 46, pushing the old row 46 (and
# everything below it) down by one row.
$ws.Rows(46).EntireRow.Insert()

# Populate the newly-inserted row with the new observation.
$ws.Range("A46").Value = 9
$ws.Range("B46").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 44468
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = 100112005
$ws.Range("G46").Value = "Puerro"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 133
$ws.Range("K46").Value = 7000
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = 7504
$ws.Range("N46").Value = "`$/paquete 20 unidades"
$ws.Range("O46").Value = "Provincia de Chacabuco"
$ws.Range("P46").Value = 375
$ws.Range("Q46").Value = 20
$ws.Range("R46").Value = "Hortaliza"
